$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper: locate a paragraph whose text contains a given substring.
# ------------------------------------------------------------------
function Find-ParagraphByText {
    param([string]$needle)
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# Paragraph "L - load program for persistent memory (auto.bas)"
#   * whole paragraph turns red
#   * a new " - TBD" run is appended (same formatting)
#   * existing runs / proofErr markers are otherwise untouched
# ------------------------------------------------------------------
$pL = Find-ParagraphByText "load program for persistent memory"
$pL.Range.Font.Color = 255

$rL = $pL.Range
$textRangeL = $d.Range($rL.Start, $rL.End - 1)

$xmlL = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve">L - </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:color w:val="FF0000"/></w:rPr><w:t>load</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> program for persistent memory (auto.bas)</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> - TBD</w:t></w:r>' +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$textRangeL.InsertXML($xmlL)

# ------------------------------------------------------------------
# Paragraph "S - save program to persistent memory (auto.bas)"
#   * whole paragraph turns red
#   * the three runs + proofErr markers collapse into a single run
#   * a new " - TBD" run is appended (same formatting)
# ------------------------------------------------------------------
$pS = Find-ParagraphByText "save program to persistent memory"
$pS.Range.Font.Color = 255

$rS = $pS.Range
$textRangeS = $d.Range($rS.Start, $rS.End - 1)

$dash = [char]0x2013
$xmlS = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
    ('<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:color w:val="FF0000"/></w:rPr><w:t>S ' + $dash + ' save program to persistent memory (auto.bas)</w:t></w:r>') +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:i/><w:color w:val="FF0000"/></w:rPr><w:t xml:space="preserve"> - TBD</w:t></w:r>' +
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$textRangeS.InsertXML($xmlS)

Write-Host "Applied red-color + TBD edits to L/S paragraphs"
